# Final patch before release: add English-language name labels in column D
# alongside the existing Japanese names (column A) and the English names
# already present for some rows in column C.
#
#   Row 2: A2 = "ライム" (Lime)   -> add D2 = "Lime"
#   Row 3: A3 = "シィナ" (Shina), C3 = "Shina" -> add D3 = "Shina"
#   Row 4: A4 = "リリー" (Lily),  C4 = "Lily"  -> add D4 = "Lily"
#
# This grows the sheet's used range from A1:C18 to A1:D18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Lime"
$ws.Range("D3").Value = "Shina"
$ws.Range("D4").Value = "Lily"
